$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 150 and 151, pushing the existing
# data (old rows 150-264) down to rows 152-266.
$ws.Rows.Item(150).Insert()
$ws.Rows.Item(151).Insert()

# Populate the first new row (150) with this week's "Primera" quality data.
$ws.Cells.Item(150, 1).Value = 1
$ws.Cells.Item(150, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(150, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(150, 4).Value = 44596
$ws.Cells.Item(150, 5).Value = 15
$ws.Cells.Item(150, 6).Value = 100112032
$ws.Cells.Item(150, 7).Value = "Zapallo italiano"
$ws.Cells.Item(150, 8).Value = "Huracán"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 170
$ws.Cells.Item(150, 11).Value = 5000
$ws.Cells.Item(150, 12).Value = 5500
$ws.Cells.Item(150, 13).Value = 5250
$ws.Cells.Item(150, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(150, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(150, 16).Value = 75
$ws.Cells.Item(150, 17).Value = 70
$ws.Cells.Item(150, 18).Value = "Hortaliza"

# Populate the second new row (151) with this week's "Segunda" quality data.
$ws.Cells.Item(151, 1).Value = 1
$ws.Cells.Item(151, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(151, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(151, 4).Value = 44596
$ws.Cells.Item(151, 5).Value = 15
$ws.Cells.Item(151, 6).Value = 100112032
$ws.Cells.Item(151, 7).Value = "Zapallo italiano"
$ws.Cells.Item(151, 8).Value = "Huracán"
$ws.Cells.Item(151, 9).Value = "Segunda"
$ws.Cells.Item(151, 10).Value = 180
$ws.Cells.Item(151, 11).Value = 4000
$ws.Cells.Item(151, 12).Value = 4500
$ws.Cells.Item(151, 13).Value = 4250
$ws.Cells.Item(151, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(151, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(151, 16).Value = 42
$ws.Cells.Item(151, 17).Value = 100
$ws.Cells.Item(151, 18).Value = "Hortaliza"
